# The source data pipeline was re-run against a renamed repo / fixed
# output-folder path, which dropped two MAGs from this species' prediction
# table (the ones whose "max" softmax score belonged to a different
# candidate species than the one this sheet is for): row 7
# (even_MAG-GUT22776.fa) and row 25 (even_MAG-GUT62790.fa). Deleting the
# higher-numbered row first keeps row 7's index valid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(25).Delete()
$ws.Rows.Item(7).Delete()
